# "updated after 2nd race"
# - add an "alternative" data source row (HradniOkruh2019-b.xlsx) to M1 and M2
# - add a new category "Ml. Zaci" to the Kategorie lookup table
# - add a brand-new "Ml. Zaci" results sheet (same layout as M1/M2) sourced
#   from the alternative file

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Kategorie" sheet: insert a new category row (row 6) for "Ml. Žáci"
# ---------------------------------------------------------------------
$wsKat = $wb.Worksheets.Item("Kategorie")
$wsKat.Rows.Item(6).Insert()
$wsKat.Range("A6").Value = "Ml. Žáci"
$wsKat.Range("B6").Value = 11
$wsKat.Range("C6").Value = 12
$wsKat.Range("D6").Formula = "=`$B`$1-C6"
$wsKat.Range("E6").Formula = "=`$B`$1-B6"

# the named range "Kategorie" grows by the inserted row
$wb.Names.Item("Kategorie").RefersTo = "=Kategorie!`$A`$4:`$E`$24"

# ---------------------------------------------------------------------
# 2) "M1" sheet: insert the alternative-source row (row 3) + new column
# ---------------------------------------------------------------------
$wsM1 = $wb.Worksheets.Item("M1")
$wsM1.Rows.Item(3).Insert()
$wsM1.Range("A3").Value = "..\2019-src\zavod1\HradniOkruh2019-b.xlsx"
$wsM1.Range("B3").Value = "M1"
$wsM1.Range("C3").Value = 2
$wsM1.Range("D3").Value = "E"
$wsM1.Range("E3").Value = "J"
$wsM1.Range("F3").Value = "I"
$wsM1.Range("G3").Value = "B"
$wsM1.Range("H3").Value = 1
$wsM1.Range("H3").Style = "Normal"
$wsM1.Range("H1").Value = "alternativní"

# ---------------------------------------------------------------------
# 3) "M2" sheet: insert the alternative-source row (row 3) + new column
# ---------------------------------------------------------------------
$wsM2 = $wb.Worksheets.Item("M2")
$wsM2.Rows.Item(3).Insert()
$wsM2.Range("A3").Value = "..\2019-src\zavod1\HradniOkruh2019-b.xlsx"
$wsM2.Range("B3").Value = "M2"
$wsM2.Range("C3").Value = 2
$wsM2.Range("D3").Value = "E"
$wsM2.Range("E3").Value = "J"
$wsM2.Range("F3").Value = "I"
$wsM2.Range("G3").Value = "B"
$wsM2.Range("H3").Value = 1
$wsM2.Range("H3").Style = "Normal"
$wsM2.Range("H1").Value = "alternativní"

# ---------------------------------------------------------------------
# 4) New "Ml. Žáci" results sheet - same layout as M1 / M2, one data row
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNew = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsNew.Name = "Ml. Žáci"

$wsM1.Range("A1:H1").Copy()
$wsNew.Range("A1:H1").PasteSpecial(-4122)

$wsNew.Range("A1").Value = "Zdroj"
$wsNew.Range("B1").Value = "Strana"
$wsNew.Range("C1").Value = "první řádek"
$wsNew.Range("D1").Value = "jmeno sl."
$wsNew.Range("E1").Value = "team sl."
$wsNew.Range("F1").Value = "rok sl."
$wsNew.Range("G1").Value = "pořadí sl."
$wsNew.Range("H1").Value = "alternativní"

$wsNew.Range("A2").Value = "..\2019-src\zavod1\HradniOkruh2019.xlsx"
$wsNew.Range("B2").Value = "D3"
$wsNew.Range("C2").Value = 2
$wsNew.Range("D2").Value = "E"
$wsNew.Range("E2").Value = "I"
$wsNew.Range("F2").Value = "H"
$wsNew.Range("G2").Value = "B"

$wsNew.Columns.Item(1).EntireColumn.AutoFit() | Out-Null

# ---------------------------------------------------------------------
# 5) Selection / active-sheet bookkeeping (mirrors the saved UI state)
# ---------------------------------------------------------------------
[void]$wsNew.Range("G2").Select()
[void]$wsM2.Range("A1:H2").Select()
[void]$wsM1.Range("F15").Select()
[void]$wsKat.Range("A7").Select()
